$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text would otherwise be auto-converted
# to a real number by Excel's input parsing; force Text format first,
# then clear the format again afterward so no stray style is left behind.
$textCells = @("D5", "D6", "D8", "D11", "D16", "D18", "D21", "D22", "D25", "D30", "D37", "D38", "D43", "D45", "D49")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated crypto price / volume figures.
$ws.Range("D2").Value = "27.654.78"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.634.01"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "212.16"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "0.524"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "23.29"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").Value = "0.0855"
$ws.Range("E11").Value = "  -4.24%  "
$ws.Range("D12").Value = "1.865.95"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "1.632.31"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "65.12"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "27.633.19"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "230.18"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "10.65"
$ws.Range("E22").Value = "  +4.22%  "
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("E24").Value = "  +3.11%  "
$ws.Range("D25").Value = "148.82"
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D33").Value = "1.481.49"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("D37").Value = "0.959"
$ws.Range("E37").Value = "  +7.13%  "
$ws.Range("D38").Value = "0.879"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +1.14%  "
$ws.Range("D43").Value = "67.65"
$ws.Range("E43").Value = "  -3.07%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "2.21"
$ws.Range("D47").Value = "1.774.60"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D49").Value = "87.54"
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("E51").Value = "  -0.20%  "

foreach ($ref in $textCells) {
    $ws.Range($ref).ClearFormats()
}
